$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '34.313.24'
$ws.Range('E2').Value = '  -0.14%  '
$ws.Range('D3').Value = '1.802.22'
$ws.Range('E3').Value = '  +0.76%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '227.65'
$ws.Range('E5').Value = '  +0.75%  '
$ws.Range('D6').Value = '0.574'
$ws.Range('E6').Value = '  +3.92%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  +11.53%  '
$ws.Range('D9').Value = '0.301'
$ws.Range('E9').Value = '  +2.24%  '
$ws.Range('D10').Value = '0.0693'
$ws.Range('E10').Value = '  +0.61%  '
$ws.Range('D11').Value = '0.0963'
$ws.Range('E11').Value = '  +2.08%  '
$ws.Range('D12').Value = '2.066.02'
$ws.Range('E12').Value = '  +0.96%  '
$ws.Range('D13').Value = '11.73'
$ws.Range('E13').Value = '  +6.50%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').Value = '0.643'
$ws.Range('E14').Value = '  +1.75%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '1.776.84'
$ws.Range('E15').Value = '  -0.80%  '
$ws.Range('D16').Value = '4.51'
$ws.Range('E16').Value = '  +5.78%  '
$ws.Range('D17').Value = '34.342.93'
$ws.Range('E17').Value = '  -0.08%  '
$ws.Range('D18').Value = '69.13'
$ws.Range('E18').Value = '  +1.46%  '
$ws.Range('D19').Value = '245.83'
$ws.Range('E19').Value = '  +0.76%  '
$ws.Range('D20').Value = '0.0₃0795'
$ws.Range('E20').Value = '  +0.09%  '
$ws.Range('D21').Value = '11.72'
$ws.Range('E21').Value = '  +4.92%  '
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('D23').Value = '4.17'
$ws.Range('E23').Value = '  +0.85%  '
$ws.Range('D24').Value = '171.71'
$ws.Range('E24').Value = '  +3.75%  '
$ws.Range('E25').Value = '  +2.79%  '
$ws.Range('D26').Value = '7.96'
$ws.Range('E26').Value = '  +9.64%  '
$ws.Range('D27').Value = '16.86'
$ws.Range('E27').Value = '  +2.46%  '
$ws.Range('D28').Value = '0.118'
$ws.Range('E28').Value = '  +2.40%  '
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('D30').Value = '4.02'
$ws.Range('E30').Value = '  +1.27%  '
$ws.Range('D31').Value = '0.0531'
$ws.Range('E31').Value = '  +1.53%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').Value = '1.24'
$ws.Range('E32').Value = '  +1.27%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '3.83'
$ws.Range('E33').Value = '  +1.11%  '
$ws.Range('D34').Value = '1.81'
$ws.Range('E34').Value = '  +0.37%  '
$ws.Range('D35').Value = '1.398.51'
$ws.Range('E35').Value = '  -0.18%  '
$ws.Range('D36').Value = '0.671'
$ws.Range('E36').Value = '  -0.54%  '
$ws.Range('D37').Value = '2.48'
$ws.Range('E37').Value = '  -4.17%  '
$ws.Range('E38').Value = '  +0.53%  '
$ws.Range('E39').Value = '  +1.07%  '
$ws.Range('B40').Value = 'WEMIXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D40').Value = '1.23'
$ws.Range('E40').Value = '  +10.64%  '
$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').Value = '0.963'
$ws.Range('E41').Value = '  +3.20%  '
$ws.Range('D42').Value = '82.39'
$ws.Range('E42').Value = '  -2.42%  '
$ws.Range('D44').Value = '2.43'
$ws.Range('E44').Value = '  +1.03%  '
$ws.Range('D45').Value = '13.40'
$ws.Range('E45').Value = '  -2.43%  '
$ws.Range('B46').Value = 'Kaspa'
$ws.Range('C46').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D46').Value = '0.0506'
$ws.Range('E46').Value = '  -3.52%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').Value = '6.03'
$ws.Range('E47').Value = '  +1.03%  '
$ws.Range('D48').Value = '1.965.77'
$ws.Range('E48').Value = '  +1.07%  '
$ws.Range('D49').Value = '104.39'
$ws.Range('E49').Value = '  -0.26%  '
$ws.Range('E50').Value = '  -0.01%  '
$ws.Range('D51').Value = '0.0₆0125'
$ws.Range('E51').Value = '  -2.08%  '
